$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 0.02023634318839907
$ws.Cells.Item(2, 7).ClearContents()
$ws.Cells.Item(2, 8).Value = -0.02915623868895263

$ws.Cells.Item(3, 6).Value = 0.1013168047700807
$ws.Cells.Item(3, 7).ClearContents()
$ws.Cells.Item(3, 8).Value = -0.03067237533543543

$ws.Cells.Item(4, 6).Value = -0.2859040023064409
$ws.Cells.Item(4, 7).ClearContents()
$ws.Cells.Item(4, 8).Value = -0.03230520711213458

$ws.Cells.Item(5, 6).Value = 0.02023634318839907
$ws.Cells.Item(5, 7).ClearContents()
$ws.Cells.Item(5, 8).Value = -0.02930843720673171

$ws.Cells.Item(6, 6).Value = -0.3767851726253837
$ws.Cells.Item(6, 7).ClearContents()
$ws.Cells.Item(6, 8).Value = -0.03145701400854761

$ws.Cells.Item(7, 6).Value = -0.1162825941674294
$ws.Cells.Item(7, 7).ClearContents()
$ws.Cells.Item(7, 8).Value = -0.04068344190870687

$ws.Cells.Item(8, 6).Value = 0.1206427006000514
$ws.Cells.Item(8, 7).ClearContents()
$ws.Cells.Item(8, 8).Value = -0.03098559925908104

$ws.Cells.Item(9, 6).Value = -0.3297628183503412
$ws.Cells.Item(9, 7).ClearContents()
$ws.Cells.Item(9, 8).Value = -0.04173023469976335

$ws.Cells.Item(10, 6).Value = 0.01233230957495361
$ws.Cells.Item(10, 7).ClearContents()
$ws.Cells.Item(10, 8).Value = -0.04013200448336178

$ws.Cells.Item(11, 6).Value = -0.3556198415185821
$ws.Cells.Item(11, 7).ClearContents()
$ws.Cells.Item(11, 8).Value = -0.04783913451522626

$ws.Cells.Item(12, 6).Value = 0.132025025699478
$ws.Cells.Item(12, 7).ClearContents()
$ws.Cells.Item(12, 8).Value = -0.04586944825788503

$ws.Cells.Item(13, 6).Value = -0.2343712103023051
$ws.Cells.Item(13, 7).ClearContents()
$ws.Cells.Item(13, 8).Value = -0.04721801424869909

$ws.Cells.Item(14, 6).Value = 0.1972432691056193
$ws.Cells.Item(14, 7).ClearContents()
$ws.Cells.Item(14, 8).Value = -0.05046050730132924

$ws.Cells.Item(15, 6).Value = 0.2956122349779051
$ws.Cells.Item(15, 7).ClearContents()
$ws.Cells.Item(15, 8).Value = -0.05410925584232507

$ws.Cells.Item(16, 6).Value = -0.21175444599907
$ws.Cells.Item(16, 7).ClearContents()
$ws.Cells.Item(16, 8).Value = -0.05860944282376654

$ws.Cells.Item(17, 6).Value = 0.08280792306821345
$ws.Cells.Item(17, 7).ClearContents()
$ws.Cells.Item(17, 8).Value = -0.05676506090363112

$ws.Cells.Item(18, 6).Value = -0.2007598915012807
$ws.Cells.Item(18, 7).ClearContents()
$ws.Cells.Item(18, 8).Value = -0.05438311111253059

$ws.Cells.Item(19, 6).Value = 0.06140240701563842
$ws.Cells.Item(19, 7).ClearContents()
$ws.Cells.Item(19, 8).Value = -0.05908256685437385

$ws.Cells.Item(20, 6).Value = 0.06140240701563842
$ws.Cells.Item(20, 7).ClearContents()
$ws.Cells.Item(20, 8).Value = -0.06373444709153075

$ws.Cells.Item(21, 6).Value = -0.2915000674038085
$ws.Cells.Item(21, 7).ClearContents()
$ws.Cells.Item(21, 8).Value = -0.06919671797134891

$ws.Cells.Item(22, 6).Value = -0.4538669545145184
$ws.Cells.Item(22, 7).ClearContents()
$ws.Cells.Item(22, 8).Value = -0.07166554474747217

$ws.Cells.Item(23, 6).Value = -0.4538669545145184
$ws.Cells.Item(23, 7).ClearContents()
$ws.Cells.Item(23, 8).Value = -0.07376537875540018

$ws.Cells.Item(24, 6).Value = -0.318272682631718
$ws.Cells.Item(24, 7).ClearContents()
$ws.Cells.Item(24, 8).Value = -0.06839602971398098

$ws.Cells.Item(25, 6).Value = -0.009077992029759957
$ws.Cells.Item(25, 7).ClearContents()
$ws.Cells.Item(25, 8).Value = -0.07383933164862536

$ws.Cells.Item(26, 6).Value = -0.4538669545145184
$ws.Cells.Item(26, 7).ClearContents()
$ws.Cells.Item(26, 8).Value = -0.09768858188808835

$ws.Cells.Item(27, 6).Value = -0.3505111980378801
$ws.Cells.Item(27, 7).ClearContents()
$ws.Cells.Item(27, 8).Value = -0.1015155303933532

$ws.Cells.Item(28, 6).Value = -0.4538669545145184
$ws.Cells.Item(28, 7).ClearContents()
$ws.Cells.Item(28, 8).Value = -0.09970411375571463

$ws.Cells.Item(29, 6).Value = 0.1987729760015876
$ws.Cells.Item(29, 7).ClearContents()
$ws.Cells.Item(29, 8).Value = -0.09642294950427645

$ws.Cells.Item(30, 6).Value = 0.1281187015163463
$ws.Cells.Item(30, 7).ClearContents()
$ws.Cells.Item(30, 8).Value = -0.07080499003679154

$ws.Cells.Item(31, 6).Value = -0.03048350808233497
$ws.Cells.Item(31, 7).ClearContents()
$ws.Cells.Item(31, 8).Value = -0.07984220631471847

$ws.Cells.Item(32, 6).Value = -0.4538669545145184
$ws.Cells.Item(32, 7).ClearContents()
$ws.Cells.Item(32, 8).Value = -0.08839472913129356
